$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F2").Value = 11
$ws.Range("F3").Value = -10
$ws.Range("F8").Value = 14
$ws.Range("F9").Value = -2
$ws.Range("F10").Value = -7
$ws.Range("F12").Value = 3
$ws.Range("F13").Value = 2
